# Fix the two "missing punch" rows (27 and 47) that were previously
# rendered with the special red/highlighted "Débito Banco Horas" style
# and flagged with empty time cells. The actual clock-in/out times were
# found, so:
#   - restore the normal (unhighlighted) style used by every other row
#   - fill in the recovered times in columns C (entrada), D (saída almoço)
#     and F (saída)
#   - remove the now-unused E/G/H/I/J/K cells (no lunch-return punch,
#     no "Débito Banco Horas" annotation needed anymore)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Row($row, $prevRow, $entrada, $saidaAlmoco, $saida) {
    # Column A (date) and B (shift marker) just need their formatting
    # reset back to the normal row style - B is already normal, A was
    # using the highlighted style.
    $ws.Range("A$row").Style = $ws.Range("A$prevRow").Style

    # Columns C, D and F get the recovered punch times, with the same
    # normal style as the surrounding rows.
    $ws.Range("C$row").Value = $entrada
    $ws.Range("C$row").Style = $ws.Range("C$prevRow").Style

    $ws.Range("D$row").Value = $saidaAlmoco
    $ws.Range("D$row").Style = $ws.Range("D$prevRow").Style

    $ws.Range("F$row").Value = $saida
    $ws.Range("F$row").Style = $ws.Range("F$prevRow").Style

    # E (lunch return), G-J (unused columns) and K ("Débito Banco Horas"
    # note) are no longer applicable, so drop them entirely.
    $ws.Range("E$row").Clear()
    $ws.Range("G$row").Clear()
    $ws.Range("H$row").Clear()
    $ws.Range("I$row").Clear()
    $ws.Range("J$row").Clear()
    $ws.Range("K$row").Clear()
}

# Row 27 - 26/07/2019 Fri
Fix-Row 27 26 "07:47" "11:46" "16:54"

# Row 47 - 15/08/2019 Thu
Fix-Row 47 46 "07:57" "11:36" "17:14"
